$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New review row (row 15) - same shape as existing rows (appid/keyword/email/
# recovery-email/time/review), reusing the com.singleton.strechy / stretchy
# keyword pairing used in row 2.
$ws.Range("A15").Value = "com.singleton.strechy"
$ws.Range("B15").Value = "stretchy"
$ws.Range("C15").Value = "itaisenior@gmail.com"
$ws.Range("D15").Value = "vikicrestina@gmail.com"
$ws.Range("E15").Value = "27/5/2019 15:59"
$ws.Range("F15").Value = "This taxi offline game is in my favourite. So much hours of playing! I love it"

# Hyperlink the email + recovery-email columns, same as every other row.
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:itaisenior@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "itaisenior@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:vikicrestina@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "vikicrestina@gmail.com")

# Match the formatting of the row above (copy it after the hyperlinks so the
# hyperlink auto-styling doesn't stomp the normal cell formats).
$ws.Range("A14:F14").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)

# Leave the selection where the author's session ended up.
$ws.Range("F15").Select()
